$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Calca"
$ws.Range("C2").Value = "Calcrl"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.354751
$ws.Range("H2").Value = 0.709502
$ws.Range("I2").Value = 0.7067029850439027
$ws.Range("J2").Value = 0.6163205031315422
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 41.6173
$ws.Range("N2").Value = 83.2346
$ws.Range("O2").Value = 0.411761355892064
$ws.Range("P2").Value = 0.3324886731607734
$ws.Range("Q2").Value = 14.7637787923
$ws.Range("R2").Value = 59.0551151692
$ws.Range("S2").Value = 0.2909929793346464
$ws.Range("T2").Value = 0.2049195863279868

$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Calcrl"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.354751
$ws.Range("H3").Value = 0.709502
$ws.Range("I3").Value = 0.7067029850439027
$ws.Range("J3").Value = 0.6163205031315422
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.980450666666668
$ws.Range("N3").Value = 26.941352
$ws.Range("O3").Value = 0.08885253351439082
$ws.Range("P3").Value = 0.1076198405427232
$ws.Range("Q3").Value = 3.185823854450667
$ws.Range("R3").Value = 19.114943126704
$ws.Range("S3").Value = 0.0627923506633334
$ws.Range("T3").Value = 0.0663283142702275

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Calca"
$ws.Range("C4").Value = "Calcrl"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.354751
$ws.Range("H4").Value = 0.709502
$ws.Range("I4").Value = 0.7067029850439027
$ws.Range("J4").Value = 0.6163205031315422
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.45113833333333
$ws.Range("N4").Value = 40.353415
$ws.Range("O4").Value = 0.1330854946963174
$ws.Range("P4").Value = 0.1611956255073737
$ws.Range("Q4").Value = 4.771804774888333
$ws.Range("R4").Value = 28.63082864933
$ws.Range("S4").Value = 0.094051916367932
$ws.Range("T4").Value = 0.09934816901530824

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Calca"
$ws.Range("C5").Value = "Calcrl"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.354751
$ws.Range("H5").Value = 0.709502
$ws.Range("I5").Value = 0.7067029850439027
$ws.Range("J5").Value = 0.6163205031315422
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.258772
$ws.Range("N5").Value = 22.517544
$ws.Range("O5").Value = 0.1113942332731726
$ws.Range("P5").Value = 0.0899485109245354
$ws.Range("Q5").Value = 3.994060625772
$ws.Range("R5").Value = 15.976242503088
$ws.Range("S5").Value = 0.07872263717082789
$ws.Range("T5").Value = 0.05543711150894268

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Calca"
$ws.Range("C6").Value = "Calcrl"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.354751
$ws.Range("H6").Value = 0.709502
$ws.Range("I6").Value = 0.7067029850439027
$ws.Range("J6").Value = 0.6163205031315422
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.450729
$ws.Range("N6").Value = 46.352187
$ws.Range("O6").Value = 0.1528694346476305
$ws.Range("P6").Value = 0.1851583014002596
$ws.Range("Q6").Value = 5.481161563479
$ws.Range("R6").Value = 32.886969380874
$ws.Range("S6").Value = 0.1080332857874543
$ws.Range("T6").Value = 0.1141168574779897

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Calca"
$ws.Range("C7").Value = "Calcrl"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.354751
$ws.Range("H7").Value = 0.709502
$ws.Range("I7").Value = 0.7067029850439027
$ws.Range("J7").Value = 0.6163205031315422
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.313018
$ws.Range("N7").Value = 30.939054
$ws.Range("O7").Value = 0.1020369479764247
$ws.Range("P7").Value = 0.1235890484643348
$ws.Range("Q7").Value = 3.658553448518
$ws.Range("R7").Value = 21.951320691108
$ws.Range("S7").Value = 0.07210981571970879
$ws.Range("T7").Value = 0.07617046453108735

$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Calca"
$ws.Range("C8").Value = "Calcrl"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1472293333333334
$ws.Range("H8").Value = 0.441688
$ws.Range("I8").Value = 0.2932970149560972
$ws.Range("J8").Value = 0.3836794968684579
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 41.6173
$ws.Range("N8").Value = 83.2346
$ws.Range("O8").Value = 0.411761355892064
$ws.Range("P8").Value = 0.3324886731607734
$ws.Range("Q8").Value = 6.127287334133334
$ws.Range("R8").Value = 36.7637240048
$ws.Range("S8").Value = 0.1207683765574176
$ws.Range("T8").Value = 0.1275690868327867

$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Calca"
$ws.Range("C9").Value = "Calcrl"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1472293333333334
$ws.Range("H9").Value = 0.441688
$ws.Range("I9").Value = 0.2932970149560972
$ws.Range("J9").Value = 0.3836794968684579
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.980450666666668
$ws.Range("N9").Value = 26.941352
$ws.Range("O9").Value = 0.08885253351439082
$ws.Range("P9").Value = 0.1076198405427232
$ws.Range("Q9").Value = 1.322185764686223
$ws.Range("R9").Value = 11.899671882176
$ws.Range("S9").Value = 0.02606018285105741
$ws.Range("T9").Value = 0.04129152627249571

$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Calca"
$ws.Range("C10").Value = "Calcrl"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1472293333333334
$ws.Range("H10").Value = 0.441688
$ws.Range("I10").Value = 0.2932970149560972
$ws.Range("J10").Value = 0.3836794968684579
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.45113833333333
$ws.Range("N10").Value = 40.353415
$ws.Range("O10").Value = 0.1330854946963174
$ws.Range("P10").Value = 0.1611956255073737
$ws.Range("Q10").Value = 1.980402129391111
$ws.Range("R10").Value = 17.82361916452
$ws.Range("S10").Value = 0.03903357832838541
$ws.Range("T10").Value = 0.06184745649206552

$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Calca"
$ws.Range("C11").Value = "Calcrl"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1472293333333334
$ws.Range("H11").Value = 0.441688
$ws.Range("I11").Value = 0.2932970149560972
$ws.Range("J11").Value = 0.3836794968684579
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 11.258772
$ws.Range("N11").Value = 22.517544
$ws.Range("O11").Value = 0.1113942332731726
$ws.Range("P11").Value = 0.0899485109245354
$ws.Range("Q11").Value = 1.657621495712
$ws.Range("R11").Value = 9.945728974272
$ws.Range("S11").Value = 0.03267159610234468
$ws.Range("T11").Value = 0.03451139941559273

$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Calca"
$ws.Range("C12").Value = "Calcrl"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1472293333333334
$ws.Range("H12").Value = 0.441688
$ws.Range("I12").Value = 0.2932970149560972
$ws.Range("J12").Value = 0.3836794968684579
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 15.450729
$ws.Range("N12").Value = 46.352187
$ws.Range("O12").Value = 0.1528694346476305
$ws.Range("P12").Value = 0.1851583014002596
$ws.Range("Q12").Value = 2.274800530184
$ws.Range("R12").Value = 20.473204771656
$ws.Range("S12").Value = 0.04483614886017622
$ws.Range("T12").Value = 0.07104144392226991

$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Calca"
$ws.Range("C13").Value = "Calcrl"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1472293333333334
$ws.Range("H13").Value = 0.441688
$ws.Range("I13").Value = 0.2932970149560972
$ws.Range("J13").Value = 0.3836794968684579
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.313018
$ws.Range("N13").Value = 30.939054
$ws.Range("O13").Value = 0.1020369479764247
$ws.Range("P13").Value = 0.1235890484643348
$ws.Range("Q13").Value = 1.518378764794667
$ws.Range("R13").Value = 13.665408883152
$ws.Range("S13").Value = 0.02992713225671596
$ws.Range("T13").Value = 0.04741858393324742
